{"js": "const body = context.document.body;\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:rPr>\n                <w:b/>\n                <w:bCs/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:bCs/>\n              </w:rPr>\n              <w:t>Technical Design Document: MarquesePendergrass_ProgrammingExercise_3.py</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:bCs/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">Name: </w:t>\n            </w:r>\n            <w:r>\n              <w:t>Marquese Pendergrass</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:bCs/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">Date Created: </w:t>\n            </w:r>\n            <w:r>\n              <w:t>September 25</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:vertAlign w:val=\"superscript\"/>\n              </w:rPr>\n              <w:t>th</w:t>\n            </w:r>\n            <w:r>\n              <w:t>, 2025</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:rPr>\n                <w:b/>\n                <w:bCs/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:bCs/>\n              </w:rPr>\n              <w:t>Program Description:</w:t>\n            </w:r>\n          </w:p>\n          <w:p/>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nbody.insertOoxml(ooxml, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$sel = $d.Content\n$sel.Collapse(1)  # wdCollapseStart\n\n$frag = @'\n<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Technical Design Document: MarquesePendergrass_ProgrammingExercise_3.py</w:t></w:r></w:p><w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\"preserve\">Name: </w:t></w:r><w:r><w:t>Marquese Pendergrass</w:t></w:r></w:p><w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=\"preserve\">Date Created: </w:t></w:r><w:r><w:t>September 25</w:t></w:r><w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>th</w:t></w:r><w:r><w:t>, 2025</w:t></w:r></w:p><w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Program Description:</w:t></w:r></w:p><w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"/>\n'@\n\n$sel.InsertXML($frag)\n"}
